# Updated cryptos list on Thu Oct  5 03:58:43 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the cryptos
# table, and re-syncs three rows (28, 49-51) whose underlying coin
# shifted rank (name/link/price/volume all changed together).
#
# Price values such as "27.696.12" / "213.03" are display text (the
# source feed uses '.' as both a thousands separator and a decimal
# point), not real numbers, so each D-column write briefly forces the
# cell to Text format before assigning the value and then clears the
# format again — this keeps the stored value a string (matching the
# original file) without leaving a lingering custom cell style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($Cell, $Text)
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-PriceText "D2" "27.696.12"
$ws.Range("E2").Value = "  +1.01%  "

# Row 3 - Ethereum
Set-PriceText "D3" "1.644.63"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
Set-PriceText "D5" "213.03"
$ws.Range("E5").Value = "  +0.45%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.85%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.09%  "

# Row 8 - Solana
Set-PriceText "D8" "23.23"
$ws.Range("E8").Value = "  +0.83%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.85%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.57%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.18%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PriceText "D12" "1.878.69"
$ws.Range("E12").Value = "  +0.19%  "

# Row 13 - WrappedEther
Set-PriceText "D13" "1.649.85"
$ws.Range("E13").Value = "  +0.43%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.26%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.21%  "

# Row 16 - Litecoin
Set-PriceText "D16" "64.81"
$ws.Range("E16").Value = "  +0.90%  "

# Row 17 - WrappedBTC
Set-PriceText "D17" "27.670.72"
$ws.Range("E17").Value = "  +1.00%  "

# Row 18 - BitcoinCash
Set-PriceText "D18" "231.29"
$ws.Range("E18").Value = "  +1.37%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +0.91%  "

# Row 20 - Chainlink
Set-PriceText "D20" "7.63"
$ws.Range("E20").Value = "  +2.32%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.13%  "

# Row 22 - Uniswap
Set-PriceText "D22" "4.29"
$ws.Range("E22").Value = "  -0.60%  "

# Row 23 - Avalanche
Set-PriceText "D23" "10.09"
$ws.Range("E23").Value = "  +8.58%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.99%  "

# Row 25 - Monero
Set-PriceText "D25" "149.98"
$ws.Range("E25").Value = "  +1.32%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -0.06%  "

# Row 27 - Stellar: unchanged

# Row 28 - now EthereumClassic (was BinanceUSD)
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-PriceText "D28" "15.68"
$ws.Range("E28").Value = "  +0.97%  "

# Row 29 - now BinanceUSD (was EthereumClassic)
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-PriceText "D29" "1.00"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.31%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +0.53%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.60%  "

# Row 33 - Maker
Set-PriceText "D33" "1.442.61"
$ws.Range("E33").Value = "  +2.34%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +1.18%  "

# Row 35 - LidoDAOToken: unchanged

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.22%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +1.41%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  +0.00%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +0.41%  "

# Row 40 - TrustWalletToken
Set-PriceText "D40" "0.883"
$ws.Range("E40").Value = "  +12.08%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -0.18%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.10%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  +2.79%  "

# Row 44 - Aave
Set-PriceText "D44" "67.09"
$ws.Range("E44").Value = "  +3.88%  "

# Row 45 - MXToken
Set-PriceText "D45" "2.24"
$ws.Range("E45").Value = "  +1.36%  "

# Row 46 - RocketPoolETH
Set-PriceText "D46" "1.788.62"
$ws.Range("E46").Value = "  +0.19%  "

# Row 47 - RenderToken
Set-PriceText "D47" "1.71"
$ws.Range("E47").Value = "  +4.54%  "

# Row 48 - Quant
Set-PriceText "D48" "85.74"
$ws.Range("E48").Value = "  -1.81%  "

# Row 49 - now Algorand (was BabyDogeCoin)
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-PriceText "D49" "0.0988"
$ws.Range("E49").Value = "  +0.20%  "

# Row 50 - now EnergySwap (was Algorand)
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceText "D50" "7.75"
$ws.Range("E50").Value = "  +1.49%  "

# Row 51 - now Cronos (was EnergySwap)
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-PriceText "D51" "0.0504"
$ws.Range("E51").Value = "  +0.93%  "
